# Update cryptocurrency price/volume figures on Sheet1 (cols D & E, rows 2-51)
# per the scraped data refresh described in the commit message / diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.807.69"
$ws.Range("E2").Value = "  -2.49%  "

$ws.Range("D3").Value = "1.775.68"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'220.72"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'0.550"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.10%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'31.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.99%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'0.0707"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.51%  "

$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").Value = "2.031.19"
$ws.Range("E12").Value = "  -0.80%  "

$ws.Range("D13").Value = "1.781.62"
$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("D14").Value = "'10.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.36%  "

$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").Value = "33.816.25"
$ws.Range("E16").Value = "  -2.25%  "

$ws.Range("D17").Value = "'4.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").Value = "'67.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.67%  "

$ws.Range("D19").Value = "'244.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.10%  "

$ws.Range("D20").Value = "0.0₃0772"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "'10.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.47%  "

$ws.Range("E23").Value = "  -3.93%  "

$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("D25").Value = "'157.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "

$ws.Range("D26").Value = "'16.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").Value = "'6.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Value = "'0.0519"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("D31").Value = "'3.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.89%  "

$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("E33").Value = "  -2.82%  "

$ws.Range("D34").Value = "'1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.25%  "

$ws.Range("D35").Value = "1.392.47"
$ws.Range("E35").Value = "  -3.62%  "

$ws.Range("D36").Value = "'0.636"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "

$ws.Range("E37").Value = "  -1.21%  "

$ws.Range("E38").Value = "  -2.31%  "

$ws.Range("D39").Value = "'0.929"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").Value = "'78.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.26%  "

$ws.Range("D42").Value = "'2.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.41%  "

$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").Value = "'5.87"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.0488"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.14%  "

$ws.Range("D46").Value = "'1.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("D47").Value = "1.928.68"
$ws.Range("E47").Value = "  -0.90%  "

$ws.Range("D48").Value = "'104.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").Value = "'11.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.70%  "

$ws.Range("D51").Value = "0.0₆0120"
$ws.Range("E51").Value = "  -2.21%  "
